# Automation of InValid login test
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# with username/password test data for an invalid-login scenario, then makes
# it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet, name it, then move it to the end of the workbook
# (after the last existing sheet).
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "InvalidLogin"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch the sheet by name (safest handle after the move) and populate
# the test data.
$ws3 = $wb.Worksheets.Item("InvalidLogin")
$ws3.Range("A1").Value = "Username"
$ws3.Range("B1").Value = "password"
$ws3.Range("A2").Value = "abcd"
$ws3.Range("B2").Value = "xyz"

# Make the new sheet the active/selected one, mirroring the author's view
# state (cell B3 selected, zoomed to 190%).
$ws3.Activate()
$ws3.Range("B3").Select()
$excel.ActiveWindow.Zoom = 190
